$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 445 (a new weekly price observation),
# pushing the existing rows 445-485 down to 446-486.
$ws.Rows.Item(445).Insert()

$ws.Range("A445").Value = 4
$ws.Range("B445").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C445").Value = "Los Lagos"
$ws.Range("D445").Value = 45106
$ws.Range("E445").Value = 10
$ws.Range("F445").Value = 100112040
$ws.Range("G445").Value = "Cilantro"
$ws.Range("H445").Value = "Sin especificar"
$ws.Range("I445").Value = "Primera"
$ws.Range("J445").Value = 60
$ws.Range("K445").Value = 14000
$ws.Range("L445").Value = 14000
$ws.Range("M445").Value = 14000
$ws.Range("N445").Value = "`$/caja 36 atados"
$ws.Range("O445").Value = "Región Metropolitana"
$ws.Range("P445").Value = 389
$ws.Range("Q445").Value = 36
$ws.Range("R445").Value = "Hortaliza"
